$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '35.428.51'
Set-TextValue $ws.Range('E2') '  +2.74%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.839.21'
Set-TextValue $ws.Range('E3') '  +1.83%  '

# Row 4
Set-TextValue $ws.Range('E4') '  +0.37%  '

# Row 5
Set-TextValue $ws.Range('D5') '231.27'
Set-TextValue $ws.Range('E5') '  +2.80%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.610'
Set-TextValue $ws.Range('E6') '  +1.54%  '

# Row 7
Set-TextValue $ws.Range('E7') '  +0.28%  '

# Row 8
Set-TextValue $ws.Range('D8') '43.82'
Set-TextValue $ws.Range('E8') '  +12.43%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.309'
Set-TextValue $ws.Range('E9') '  +7.09%  '

# Row 10
Set-TextValue $ws.Range('E10') '  +5.32%  '

# Row 11
Set-TextValue $ws.Range('E11') '  +2.65%  '

# Row 12
Set-TextValue $ws.Range('D12') '2.103.87'

# Row 13
Set-TextValue $ws.Range('D13') '1.842.50'
Set-TextValue $ws.Range('E13') '  +2.05%  '

# Row 14
Set-TextValue $ws.Range('D14') '11.27'
Set-TextValue $ws.Range('E14') '  +1.63%  '

# Row 15
Set-TextValue $ws.Range('E15') '  +6.67%  '

# Row 16
Set-TextValue $ws.Range('D16') '4.70'
Set-TextValue $ws.Range('E16') '  +7.16%  '

# Row 17
Set-TextValue $ws.Range('D17') '35.377.91'
Set-TextValue $ws.Range('E17') '  +2.67%  '

# Row 18
Set-TextValue $ws.Range('D18') '70.07'
Set-TextValue $ws.Range('E18') '  +2.79%  '

# Row 19
Set-TextValue $ws.Range('D19') '0.0₃0802'
Set-TextValue $ws.Range('E19') '  +4.20%  '

# Row 20
Set-TextValue $ws.Range('D20') '244.48'
Set-TextValue $ws.Range('E20') '  +1.35%  '

# Row 21
Set-TextValue $ws.Range('D21') '12.11'
Set-TextValue $ws.Range('E21') '  +7.85%  '

# Row 22
Set-TextValue $ws.Range('E22') '  +14.40%  '

# Row 23
Set-TextValue $ws.Range('E23') '  +0.25%  '

# Row 24
Set-TextValue $ws.Range('E24') '  +1.07%  '

# Row 25
Set-TextValue $ws.Range('D25') '170.60'
Set-TextValue $ws.Range('E25') '  -0.21%  '

# Row 26
Set-TextValue $ws.Range('D26') '7.91'
Set-TextValue $ws.Range('E26') '  +2.95%  '

# Row 27
Set-TextValue $ws.Range('E27') '  +0.89%  '

# Row 28
Set-TextValue $ws.Range('E28') '  -0.69%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.54'
Set-TextValue $ws.Range('E29') '  +26.09%  '

# Row 30
Set-TextValue $ws.Range('E30') '  +0.42%  '

# Row 31
Set-TextValue $ws.Range('D31') '3.342.97'
Set-TextValue $ws.Range('E31') '  +37.59%  '

# Row 32
Set-TextValue $ws.Range('D32') '0.0554'
Set-TextValue $ws.Range('E32') '  +7.80%  '

# Row 33
Set-TextValue $ws.Range('D33') '4.09'
Set-TextValue $ws.Range('E33') '  +6.25%  '

# Row 34
Set-TextValue $ws.Range('E34') '  +4.49%  '

# Row 35
Set-TextValue $ws.Range('D35') '1.85'
Set-TextValue $ws.Range('E35') '  +1.63%  '

# Row 36
Set-TextValue $ws.Range('D36') '95.96'
Set-TextValue $ws.Range('E36') '  +16.27%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.687'
Set-TextValue $ws.Range('E37') '  +7.45%  '

# Row 38
Set-TextValue $ws.Range('B38') 'TrustWalletToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D38') '1.11'
Set-TextValue $ws.Range('E38') '  +5.33%  '

# Row 39
Set-TextValue $ws.Range('B39') 'InjectiveProtocol'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D39') '15.64'
Set-TextValue $ws.Range('E39') '  +11.95%  '

# Row 40
Set-TextValue $ws.Range('D40') '1.346.36'
Set-TextValue $ws.Range('E40') '  +3.10%  '

# Row 41
Set-TextValue $ws.Range('B41') 'RenderToken'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D41') '2.44'
Set-TextValue $ws.Range('E41') '  +5.02%  '

# Row 42
Set-TextValue $ws.Range('B42') 'VeChain'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D42') '0.0195'
Set-TextValue $ws.Range('E42') '  +4.74%  '

# Row 43
Set-TextValue $ws.Range('E43') '  +5.92%  '

# Row 44
Set-TextValue $ws.Range('E44') '  +4.62%  '

# Row 46
Set-TextValue $ws.Range('D46') '2.82'
Set-TextValue $ws.Range('E46') '  +0.65%  '

# Row 47
Set-TextValue $ws.Range('D47') '6.27'
Set-TextValue $ws.Range('E47') '  +8.67%  '

# Row 48
Set-TextValue $ws.Range('E48') '  +1.28%  '

# Row 49
Set-TextValue $ws.Range('D49') '2.006.50'
Set-TextValue $ws.Range('E49') '  +1.99%  '

# Row 50
Set-TextValue $ws.Range('E50') '  +0.31%  '

# Row 51
Set-TextValue $ws.Range('D51') '103.05'
Set-TextValue $ws.Range('E51') '  +0.18%  '
